# Updated cryptos list on Tue Jul  2 08:10:05 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.572.70'
$ws.Range('D3').Value = '3.436.33'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.43'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.37'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.480'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.95'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +3.24%  '
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('E11').Value = '  +2.06%  '
$ws.Range('D12').Value = '4.026.41'
$ws.Range('E12').Value = '  -1.49%  '
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.16'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -5.52%  '
$ws.Range('D15').Value = '3.441.40'
$ws.Range('E15').Value = '  -1.34%  '
$ws.Range('E16').Value = '  -1.48%  '
$ws.Range('D17').Value = '62.681.27'
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.56'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +1.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.03'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -3.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '386.64'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.21'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.560'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -1.21%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '3.583.14'
$ws.Range('E25').Value = '  -1.31%  '
$ws.Range('E26').Value = '  -3.17%  '
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.59'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -2.43%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.96'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -4.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.11'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -1.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('E33').Value = '  -7.66%  '
$ws.Range('E34').Value = '  -2.99%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.32'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -1.07%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.62'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +2.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '31.80'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.96'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -2.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '170.58'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').Value = '3.472.51'
$ws.Range('E40').Value = '  -1.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0771'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.49'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.69'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -2.33%  '
$ws.Range('E45').Value = '  -3.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.18'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -3.26%  '
$ws.Range('D47').Value = '2.560.75'
$ws.Range('E47').Value = '  -2.73%  '
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.25'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -2.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.69'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -3.92%  '
$ws.Range('E51').Value = '  +0.00%  '
